# Slide 8 ("Collaborations") - Content Placeholder bullet list update:
#  - Rename "Cross time-series operators:" -> "Graph-TS cross-operators:" (now the
#    first occurrence, replacing/merging with the old duplicate heading)
#  - Move the "Correlate graph metrics..." / "e.g., landslide monitoring..." /
#    "Correlate soil drying..." / "Multistore:" block up, right after
#    "Identify plants/grids..."
#  - Drop the now-redundant duplicate "Graph-TS cross-operators:" paragraph
#  - Append a trailing period to the final "Provide a unified language..." bullet

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$lines = @(
    "Physical level:",
    "TS data require different data layout than graph data;",
    "LSM-Tree-like (e.g., RocksDB)",
    "InfluxDB 3.0 on Parquet.",
    "Metadata modelling (ausiliary structures);",
    "query formalization and optimization;",
    "Analytics:",
    "TS operators in Cypher/GQL (Graph analytics);",
    "shape/patthern matching;",
    "Graph-TS cross-operators:",
    "Identify plants/grids with similar drying patterns over the last 24h",
    "Correlate graph metrics with time-series trends",
    "e.g., landslide monitoring sensor network: correlation between pressure measurements and node degree between nearby sensors",
    "Correlate soil drying with temperature (spatial join with ARPAE weather stations)",
    "Multistore:",
    "Provide a unified language that transparently distributes the execution plan on different engines."
)

$levels = @(1,2,3,3,2,2,1,2,3,2,3,3,4,3,1,2)

$tr.Text = [string]::Join([char]13, $lines)

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $tr.Paragraphs($i, 1).IndentLevel = $levels[$i-1]
}
